$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextCell $ws 'D2' '45.009.38'
Set-TextCell $ws 'E2' '  +0.62%  '

Set-TextCell $ws 'D3' '2.266.53'
Set-TextCell $ws 'E3' '  +0.67%  '

Set-TextCell $ws 'E4' '  -0.65%  '

Set-TextCell $ws 'D5' '302.01'
Set-TextCell $ws 'E5' '  -1.84%  '

Set-TextCell $ws 'D6' '94.52'
Set-TextCell $ws 'E6' '  -1.60%  '

Set-TextCell $ws 'D7' '0.565'
Set-TextCell $ws 'E7' '  -1.37%  '

Set-TextCell $ws 'D8' '1.00'
Set-TextCell $ws 'E8' '  -0.52%  '

Set-TextCell $ws 'D9' '0.509'
Set-TextCell $ws 'E9' '  -2.67%  '

Set-TextCell $ws 'D10' '34.32'
Set-TextCell $ws 'E10' '  -3.49%  '

Set-TextCell $ws 'D11' '0.0788'
Set-TextCell $ws 'E11' '  -2.36%  '

Set-TextCell $ws 'D12' '7.21'
Set-TextCell $ws 'E12' '  -0.81%  '

Set-TextCell $ws 'E13' '  -0.79%  '

Set-TextCell $ws 'D14' '2.614.04'
Set-TextCell $ws 'E14' '  +0.67%  '

Set-TextCell $ws 'D15' '2.271.12'
Set-TextCell $ws 'E15' '  -1.50%  '

Set-TextCell $ws 'D16' '13.74'
Set-TextCell $ws 'E16' '  +0.75%  '

Set-TextCell $ws 'D17' '0.798'
Set-TextCell $ws 'E17' '  -5.13%  '

Set-TextCell $ws 'D18' '44.881.51'
Set-TextCell $ws 'E18' '  +0.84%  '

Set-TextCell $ws 'D19' '12.88'
Set-TextCell $ws 'E19' '  +6.76%  '

Set-TextCell $ws 'D20' '0.0₃0922'
Set-TextCell $ws 'E20' '  -3.68%  '

Set-TextCell $ws 'D21' '6.09'
Set-TextCell $ws 'E21' '  -3.75%  '

Set-TextCell $ws 'D22' '65.38'
Set-TextCell $ws 'E22' '  -0.51%  '

Set-TextCell $ws 'D23' '239.06'
Set-TextCell $ws 'E23' '  +0.09%  '

Set-TextCell $ws 'E24' '  -3.08%  '

Set-TextCell $ws 'D25' '0.997'
Set-TextCell $ws 'E25' '  -0.54%  '

Set-TextCell $ws 'D26' '1.92'
Set-TextCell $ws 'E26' '  -4.96%  '

Set-TextCell $ws 'D27' '41.66'
Set-TextCell $ws 'E27' '  +10.72%  '

Set-TextCell $ws 'E28' '  +0.55%  '

Set-TextCell $ws 'D29' '9.55'
Set-TextCell $ws 'E29' '  -3.16%  '

Set-TextCell $ws 'D30' '19.54'
Set-TextCell $ws 'E30' '  -2.35%  '

Set-TextCell $ws 'D31' '151.92'
Set-TextCell $ws 'E31' '  -0.50%  '

Set-TextCell $ws 'E32' '  -7.10%  '

Set-TextCell $ws 'D33' '0.0787'
Set-TextCell $ws 'E33' '  -1.62%  '

Set-TextCell $ws 'D34' '2.58'
Set-TextCell $ws 'E34' '  -2.45%  '

Set-TextCell $ws 'D35' '2.94'
Set-TextCell $ws 'E35' '  -4.18%  '

Set-TextCell $ws 'E36' '  -1.77%  '

Set-TextCell $ws 'D37' '0.105'
Set-TextCell $ws 'E37' '  -4.54%  '

Set-TextCell $ws 'D38' '1.75'
Set-TextCell $ws 'E38' '  -5.94%  '

Set-TextCell $ws 'D39' '3.86'
Set-TextCell $ws 'E39' '  +1.63%  '

Set-TextCell $ws 'E40' '  +1.47%  '

Set-TextCell $ws 'D41' '3.26'
Set-TextCell $ws 'E41' '  -5.08%  '

Set-TextCell $ws 'D42' '13.83'
Set-TextCell $ws 'E42' '  -8.07%  '

Set-TextCell $ws 'D43' '1.00'
Set-TextCell $ws 'E43' '  -0.75%  '

Set-TextCell $ws 'D44' '1.95'
Set-TextCell $ws 'E44' '  +11.78%  '

Set-TextCell $ws 'D45' '1.766.03'
Set-TextCell $ws 'E45' '  -3.55%  '

Set-TextCell $ws 'D46' '0.194'
Set-TextCell $ws 'E46' '  +0.25%  '

Set-TextCell $ws 'D47' '70.30'
Set-TextCell $ws 'E47' '  -1.05%  '

Set-TextCell $ws 'B48' 'Aave'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D48' '96.84'
Set-TextCell $ws 'E48' '  -2.76%  '

Set-TextCell $ws 'B49' 'BitcoinSV'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextCell $ws 'D49' '75.44'
Set-TextCell $ws 'E49' '  -5.36%  '

Set-TextCell $ws 'D50' '7.86'
Set-TextCell $ws 'E50' '  -2.65%  '

Set-TextCell $ws 'D51' '53.13'
Set-TextCell $ws 'E51' '  -3.47%  '
